$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the footer text box that holds the
# "... Package version  0.5.5 •  Updated: 2021-07" line (shape id 322).
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame -and $cand.TextFrame.HasText) {
        $t = $cand.TextFrame.TextRange.Text
        if ($t.IndexOf("Package version") -ge 0) {
            $sh = $cand
            break
        }
    }
}

if ($sh -ne $null) {
    $tr = $sh.TextFrame.TextRange

    # Bump the package version digit: "0.5.5" -> "0.5.7"
    $full = $tr.Text
    $verIdx = $full.IndexOf("0.5.5")
    if ($verIdx -ge 0) {
        $pos = $verIdx + 5
        $tr.Characters($pos, 1).Text = "7"
    }

    # Bump the "Updated" month: "2021-07" -> "2021-08"
    $full = $tr.Text
    $dateIdx = $full.IndexOf("2021-07")
    if ($dateIdx -ge 0) {
        $pos = $dateIdx + 6
        $tr.Characters($pos, 2).Text = "08"
    }
}
